$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1442.6296
$ws.Range("I98").Value = 1218.04
$ws.Range("K98").Value = 1218.04
$ws.Range("M98").Value = 279.96
$ws.Range("H122").Value = 1442.6296
$ws.Range("I122").Value = 1218.04
$ws.Range("K122").Value = 3654.12
$ws.Range("M122").Value = -1204.12
$ws.Range("H132").Value = 1043.2195
$ws.Range("I132").Value = 975.4211
$ws.Range("K132").Value = 2926.2633
$ws.Range("M132").Value = -396.2633000000001
$ws.Range("H138").Value = 3328.836
$ws.Range("I138").Value = 4320.4736
$ws.Range("J138").Value = 2936.3125
$ws.Range("K138").Value = 12961.4208
$ws.Range("L138").Value = 8808.9375
$ws.Range("M138").Value = -7821.4208
$ws.Range("N138").Value = -19088.9375

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 214468.11
$ws.Range("I2").Value = 278388.6
$ws.Range("K2").Value = 278388.6
$ws.Range("M2").Value = -278275.6
$ws.Range("H32").Value = 18871.127
$ws.Range("I32").Value = 16081.55
$ws.Range("J32").Value = 20937.482
$ws.Range("K32").Value = 16081.55
$ws.Range("L32").Value = 20937.482
$ws.Range("M32").Value = -15794.55
$ws.Range("N32").Value = -21511.482
$ws.Range("H41").Value = 27529.5
$ws.Range("I41").Value = 17685.334
$ws.Range("K41").Value = 17685.334
$ws.Range("M41").Value = -17271.334
$ws.Range("H45").Value = 1334.3158
$ws.Range("J45").Value = 1957
$ws.Range("L45").Value = 1957
$ws.Range("N45").Value = -2711
$ws.Range("H61").Value = 31258.2
$ws.Range("I61").Value = 55734.8
$ws.Range("J61").Value = 6781.6
$ws.Range("K61").Value = 55734.8
$ws.Range("L61").Value = 6781.6
$ws.Range("M61").Value = -55522.8
$ws.Range("N61").Value = -7205.6
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H116").Value = 214468.11
$ws.Range("I116").Value = 278388.6
$ws.Range("K116").Value = 278388.6
$ws.Range("M116").Value = -276094.6
$ws.Range("H132").Value = 3531.5334
$ws.Range("I132").Value = 4418.75
$ws.Range("J132").Value = 3208.9092
$ws.Range("K132").Value = 13256.25
$ws.Range("L132").Value = 9626.7276
$ws.Range("M132").Value = -10726.25
$ws.Range("N132").Value = -14686.7276
$ws.Range("H136").Value = 31258.2
$ws.Range("I136").Value = 55734.8
$ws.Range("J136").Value = 6781.6
$ws.Range("K136").Value = 167204.4
$ws.Range("L136").Value = 20344.8
$ws.Range("M136").Value = -164654.4
$ws.Range("N136").Value = -25444.8

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 214468.11
$ws.Range("I3").Value = 278388.6
$ws.Range("K3").Value = 278388.6
$ws.Range("M3").Value = -278274.6
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H94").Value = 1532.2858
$ws.Range("J94").Value = 1844.2
$ws.Range("L94").Value = 1844.2
$ws.Range("N94").Value = -2746.2
$ws.Range("H100").Value = 54895
$ws.Range("J100").Value = 54895
$ws.Range("L100").Value = 54895
$ws.Range("N100").Value = -57059
$ws.Range("H134").Value = 5214.1514
$ws.Range("I134").Value = 5393.9614
$ws.Range("K134").Value = 16181.8842
$ws.Range("M134").Value = -13646.8842

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3107.2415
$ws.Range("I31").Value = 1369.6875
$ws.Range("J31").Value = 5245.769
$ws.Range("K31").Value = 1369.6875
$ws.Range("L31").Value = 5245.769
$ws.Range("M31").Value = -1074.6875
$ws.Range("N31").Value = -5835.769
$ws.Range("H34").Value = 3107.2415
$ws.Range("I34").Value = 1369.6875
$ws.Range("J34").Value = 5245.769
$ws.Range("K34").Value = 1369.6875
$ws.Range("L34").Value = 5245.769
$ws.Range("M34").Value = -1167.6875
$ws.Range("N34").Value = -5649.769
$ws.Range("H74").Value = 34999.2
$ws.Range("J74").Value = 34999.2
$ws.Range("L74").Value = 34999.2
$ws.Range("N74").Value = -36747.2
$ws.Range("H77").Value = 34999.2
$ws.Range("J77").Value = 34999.2
$ws.Range("L77").Value = 104997.6
$ws.Range("N77").Value = -113733.6
$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992
$ws.Range("H96").Value = 25997.5
$ws.Range("J96").Value = 25997.5
$ws.Range("L96").Value = 25997.5
$ws.Range("N96").Value = -31489.5
$ws.Range("H106").Value = 19000
$ws.Range("J106").Value = 19000
$ws.Range("L106").Value = 19000
$ws.Range("N106").Value = -21524
$ws.Range("H132").Value = 3251.2
$ws.Range("I132").Value = 2128.5
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 6385.5
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -3855.5
$ws.Range("N132").Value = -17059.0001
$ws.Range("H141").Value = 67510.82000000001
$ws.Range("J141").Value = 67061.89999999999
$ws.Range("L141").Value = 67061.89999999999
$ws.Range("N141").Value = -77421.89999999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 26052.568
$ws.Range("J131").Value = 28307.941
$ws.Range("L131").Value = 84923.823
$ws.Range("N131").Value = -95003.823

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H122").Value = 1968
$ws.Range("I122").Value = 1776.1666
$ws.Range("J122").Value = 2351.6667
$ws.Range("K122").Value = 5328.4998
$ws.Range("L122").Value = 7055.000100000001
$ws.Range("M122").Value = -2878.4998
$ws.Range("N122").Value = -11955.0001
$ws.Range("H132").Value = 1485255.5
$ws.Range("I132").Value = 1837879.5
$ws.Range("K132").Value = 5513638.5
$ws.Range("M132").Value = -5511108.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H104").Value = 21499
$ws.Range("J104").Value = 21499
$ws.Range("L104").Value = 21499
$ws.Range("N104").Value = -28487
$ws.Range("H132").Value = 3738.561
$ws.Range("I132").Value = 3159.9375
$ws.Range("J132").Value = 4108.88
$ws.Range("K132").Value = 9479.8125
$ws.Range("L132").Value = 12326.64
$ws.Range("M132").Value = -6949.8125
$ws.Range("N132").Value = -17386.64
$ws.Range("H136").Value = 1763.2142
$ws.Range("I136").Value = 1568.9
$ws.Range("J136").Value = 2249
$ws.Range("K136").Value = 4706.700000000001
$ws.Range("L136").Value = 6747
$ws.Range("M136").Value = -2156.700000000001
$ws.Range("N136").Value = -11847

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1824.9
$ws.Range("I81").Value = 916.55554
$ws.Range("K81").Value = 1833.11108
$ws.Range("M81").Value = -772.1110799999999
$ws.Range("H84").Value = 1824.9
$ws.Range("I84").Value = 916.55554
$ws.Range("K84").Value = 9165.555399999999
$ws.Range("M84").Value = -3861.555399999999
$ws.Range("H101").Value = 11148.857
$ws.Range("J101").Value = 11148.857
$ws.Range("L101").Value = 11148.857
$ws.Range("N101").Value = -17638.857
$ws.Range("H132").Value = 1510.925
$ws.Range("I132").Value = 1277.8438
$ws.Range("J132").Value = 2443.25
$ws.Range("K132").Value = 3833.5314
$ws.Range("L132").Value = 7329.75
$ws.Range("M132").Value = -1303.5314
$ws.Range("N132").Value = -12389.75
